# Add test cases for space-less SQL injections
# Target sheet: "UsernameValidatorTestData" (sheet1) in TestData.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UsernameValidatorTestData")

# ---------------------------------------------------------------------------
# 1. Remove the old row 13 ("INSERT " SQL-injection case) -- this shifts the
#    remaining rows (14-17) up by one, carrying their row height / style
#    along (old row 14, with its 30pt height and wrap-text cell, becomes the
#    new row 13).
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2. Append four brand-new test rows (17-20) with space-less SQL injection
#    strings for the username field.
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = "SELECT/*avoid-spaces*/password/**/FROM/**/Members"
$ws.Range("B18").Value = "select'asdf'as[asdf]into[#MyTable]"
$ws.Range("B19").Value = "DELETE/**/FROM/**/Users"
$ws.Range("B20").Value = "exec[sp_executesql]N'select''asdf''as[asdf]into[#MyTable]'"

$ws.Range("A17").Value = 16
$ws.Range("C17").Value = "INVALID"
$ws.Range("D17").Value = "Doesn't contain invalid characters."

$ws.Range("A18").Value = 17
$ws.Range("C18").Value = "INVALID"
$ws.Range("D18").Value = "Doesn't contain invalid characters."

$ws.Range("A19").Value = 18
$ws.Range("C19").Value = "INVALID"
$ws.Range("D19").Value = "Doesn't contain invalid characters."

$ws.Range("A20").Value = 19
$ws.Range("C20").Value = "INVALID"
$ws.Range("D20").Value = "Doesn't contain invalid characters."

# ---------------------------------------------------------------------------
# 3. Flatten every formula in A2:D16 down to its cached static value (a
#    "Paste Special -> Values" style cleanup), row by row / column by
#    column, matching the data that was already being computed.
# ---------------------------------------------------------------------------
$ws.Range("D6").Value = "Contains at least 4 characters."
$ws.Range("D7").Value = "Contains at least 4 characters."

for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 3).Value = "VALID"
}
for ($r = 6; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = "INVALID"
}

for ($r = 10; $r -le 16; $r++) {
    $ws.Cells.Item($r, 4).Value = "*"
}

# ---------------------------------------------------------------------------
# 4. Drop the now non-existent test input for row 10 (blank "Test Input"
#    cell, previously just carrying a style with no content).
# ---------------------------------------------------------------------------
$ws.Range("B10").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 5. Strip the inherited text-format styling (numFmt 49, "@") from every
#    data cell -- after the values are flattened the cells revert to
#    General formatting.
# ---------------------------------------------------------------------------
$ws.Range("A1:D16").Style = "Normal"
$ws.Range("A17:D20").Style = "Normal"

# ---------------------------------------------------------------------------
# 6. Re-apply wrap text to the long multi-line SQL-injection cell (now at
#    B13, used to be B14).
# ---------------------------------------------------------------------------
$ws.Range("B13").Style = "Normal"
$ws.Range("B13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 30

# ---------------------------------------------------------------------------
# 7. Sheet cosmetics: wider "Test Input" column, updated selection &
#    dimension follow automatically from the writes above, but make sure
#    the selection matches what was saved.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 93.85546875
$ws.Range("A1:D20").Select() | Out-Null
